$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1653
$ws.Range("B3").Value = 1336
$ws.Range("B4").Value = 1322
$ws.Range("B5").Value = 1412
$ws.Range("B6").Value = 1298
$ws.Range("B7").Value = 1264
$ws.Range("B8").Value = 1499
$ws.Range("C8").Value = 134
$ws.Range("B9").Value = 1110
$ws.Range("C9").Value = 137
$ws.Range("B10").Value = 1086
$ws.Range("C10").Value = 144
$ws.Range("B11").Value = 816
$ws.Range("C11").Value = 166
$ws.Range("B12").Value = 588
$ws.Range("C12").Value = 181
$ws.Range("B13").Value = 1622
$ws.Range("B14").Value = 1503
$ws.Range("B15").Value = 1654
$ws.Range("B16").Value = 1149
$ws.Range("B17").Value = 1536
$ws.Range("B18").Value = 1267
$ws.Range("B19").Value = 1431
$ws.Range("C19").Value = 138
$ws.Range("B20").Value = 1338
$ws.Range("C20").Value = 145
$ws.Range("B21").Value = 1265
$ws.Range("C21").Value = 158
$ws.Range("B22").Value = 1143
$ws.Range("C22").Value = 174
$ws.Range("B23").Value = 947
$ws.Range("C23").Value = 163
$ws.Range("B24").Value = 1587
$ws.Range("B25").Value = 1158
$ws.Range("B26").Value = 1707
$ws.Range("B27").Value = 1113
$ws.Range("B28").Value = 1327
$ws.Range("B29").Value = 1548
$ws.Range("B30").Value = 1299
$ws.Range("C30").Value = 146
$ws.Range("B31").Value = 1185
$ws.Range("C31").Value = 134
$ws.Range("B32").Value = 1103
$ws.Range("C32").Value = 149
$ws.Range("B33").Value = 968
$ws.Range("C33").Value = 156
$ws.Range("B34").Value = 740
$ws.Range("C34").Value = 171
$ws.Range("B35").Value = 1265
$ws.Range("B36").Value = 1521
$ws.Range("B37").Value = 1341
$ws.Range("B38").Value = 1614
$ws.Range("B39").Value = 1804
$ws.Range("B40").Value = 1462
$ws.Range("B41").Value = 1210
$ws.Range("C41").Value = 130
$ws.Range("B42").Value = 1367
$ws.Range("C42").Value = 132
$ws.Range("B43").Value = 1112
$ws.Range("C43").Value = 147
$ws.Range("B44").Value = 943
$ws.Range("C44").Value = 185
$ws.Range("B45").Value = 557
$ws.Range("C45").Value = 166
$ws.Range("B46").Value = 1669
$ws.Range("B47").Value = 1576
$ws.Range("B48").Value = 1672
$ws.Range("B49").Value = 1440
$ws.Range("B50").Value = 1133
$ws.Range("B51").Value = 1613
$ws.Range("B52").Value = 1193
$ws.Range("C52").Value = 128
$ws.Range("B53").Value = 1237
$ws.Range("C53").Value = 161
$ws.Range("B54").Value = 1063
$ws.Range("C54").Value = 164
$ws.Range("B55").Value = 803
$ws.Range("C55").Value = 178
$ws.Range("B56").Value = 610
$ws.Range("C56").Value = 195
$ws.Range("B57").Value = 1406
$ws.Range("B58").Value = 1525
$ws.Range("B59").Value = 1494
$ws.Range("B60").Value = 1348
$ws.Range("B61").Value = 1251
$ws.Range("B62").Value = 1238
$ws.Range("B63").Value = 1178
$ws.Range("C63").Value = 116
$ws.Range("B64").Value = 1147
$ws.Range("C64").Value = 131
$ws.Range("B65").Value = 1007
$ws.Range("C65").Value = 169
$ws.Range("B66").Value = 891
$ws.Range("C66").Value = 143
$ws.Range("B67").Value = 632
$ws.Range("C67").Value = 183
$ws.Range("B68").Value = 1493
$ws.Range("B69").Value = 1146
$ws.Range("B70").Value = 1323
$ws.Range("B71").Value = 1623
$ws.Range("B72").Value = 1042
$ws.Range("B73").Value = 1374
$ws.Range("B74").Value = 1338
$ws.Range("C74").Value = 134
$ws.Range("B75").Value = 1051
$ws.Range("C75").Value = 121
$ws.Range("B76").Value = 1010
$ws.Range("C76").Value = 160
$ws.Range("B77").Value = 864
$ws.Range("C77").Value = 145
$ws.Range("B78").Value = 584
$ws.Range("C78").Value = 178
$ws.Range("B79").Value = 1892
$ws.Range("B80").Value = 1318
$ws.Range("B81").Value = 1309
$ws.Range("B82").Value = 1283
$ws.Range("B83").Value = 1376
$ws.Range("B84").Value = 1204
$ws.Range("B85").Value = 1132
$ws.Range("C85").Value = 116
$ws.Range("B86").Value = 1344
$ws.Range("C86").Value = 159
$ws.Range("B87").Value = 1088
$ws.Range("C87").Value = 138
$ws.Range("B88").Value = 919
$ws.Range("C88").Value = 165
$ws.Range("B89").Value = 595
$ws.Range("C89").Value = 187
